$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 2.55
$ws.Range("I2").Value = 2.45
$ws.Range("J2").Value = 3.2
$ws.Range("Q2").Value = 1.65
$ws.Range("R2").Value = 2.2
$ws.Range("AI2").Value = 13
$ws.Range("AJ2").Value = 9.5
$ws.Range("AK2").Value = 23

# Row 3 updates
$ws.Range("Q3").Value = 2.6
$ws.Range("R3").Value = 1.48
